$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row height changes (rows 13 and 14 shrink)
$ws.Rows.Item(13).RowHeight = 68
$ws.Rows.Item(14).RowHeight = 85

# New row 17: Questionnaire / "questionnaire" Task.input
$ws.Range("H17").Value = "<span class=""bg-success"" markdown=""1"">“questionnaire” Task.input </span><!-- new-content -->"
$ws.Range("B17").Value = "<span class=""bg-success"" markdown=""1"">Questionnaire</span><!-- new-content -->"
$ws.Range("A17").Value = 16

# New row 18: QuestionnaireResponse / Attachment.Content
$ws.Range("B18").Value = "<span class=""bg-success"" markdown=""1"">QuestionnaireResponse</span><!-- new-content -->"
$ws.Range("G18").Value = "<span class=""bg-success"" markdown=""1"">Attachment.Content</span><!-- new-content -->"
$ws.Range("A18").Value = 17

# Restore selection to match the saved view
$ws.Range("F24").Select()
